$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed values (rows 2-6) ---
$ws.Range("D2").Value = 11517
$ws.Range("E2").Value = -261
$ws.Range("F2").Value = -261
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 2
$ws.Range("K2").Value = 9632
$ws.Range("L2").Value = 1393
$ws.Range("M2").Value = 8239
$ws.Range("N2").Value = 8239
$ws.Range("P2").Value = 44
$ws.Range("Q2").Value = -498
$ws.Range("R2").Value = 477
$ws.Range("S2").Value = -9
$ws.Range("T2").Value = 302
$ws.Range("U2").Value = -800
$ws.Range("W2").Value = -2.26
$ws.Range("X2").Value = 0.01
$ws.Range("Y2").Value = 0.02
$ws.Range("Z2").Value = 0.02
$ws.Range("AA2").Value = 16.91
$ws.Range("AB2").Value = 19005.86
$ws.Range("AC2").Value = 185
$ws.Range("AD2").Value = 3448.54
$ws.Range("AE2").Value = 973486
$ws.Range("AF2").Value = 0.66
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 0.16
$ws.Range("AI2").Value = 520.23
$ws.Range("AJ2").Value = 720000
$ws.Range("D3").Value = 12150
$ws.Range("E3").Value = 201
$ws.Range("F3").Value = 201
$ws.Range("G3").Value = 380
$ws.Range("H3").Value = 267
$ws.Range("I3").Value = 267
$ws.Range("K3").Value = 9941
$ws.Range("L3").Value = 1450
$ws.Range("M3").Value = 8491
$ws.Range("N3").Value = 8491
$ws.Range("P3").Value = 44
$ws.Range("Q3").Value = 553
$ws.Range("R3").Value = -601
$ws.Range("S3").Value = 7
$ws.Range("T3").Value = 257
$ws.Range("U3").Value = 296
$ws.Range("W3").Value = 1.66
$ws.Range("X3").Value = 2.19
$ws.Range("Y3").Value = 3.19
$ws.Range("Z3").Value = 2.72
$ws.Range("AA3").Value = 17.08
$ws.Range("AB3").Value = 19572.94
$ws.Range("AC3").Value = 30064
$ws.Range("AD3").Value = 24.55
$ws.Range("AE3").Value = 1003278
$ws.Range("AF3").Value = 0.74
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 0.14
$ws.Range("AI3").Value = 3.21
$ws.Range("AJ3").Value = 720000
$ws.Range("D4").Value = 12392
$ws.Range("E4").Value = 418
$ws.Range("F4").Value = 418
$ws.Range("G4").Value = 460
$ws.Range("H4").Value = 372
$ws.Range("I4").Value = 372
$ws.Range("K4").Value = 10415
$ws.Range("L4").Value = 1522
$ws.Range("M4").Value = 8893
$ws.Range("N4").Value = 8893
$ws.Range("P4").Value = 44
$ws.Range("Q4").Value = 1210
$ws.Range("R4").Value = -783
$ws.Range("S4").Value = -9
$ws.Range("T4").Value = 392
$ws.Range("U4").Value = 818
$ws.Range("W4").Value = 3.38
$ws.Range("X4").Value = 3
$ws.Range("Y4").Value = 4.28
$ws.Range("Z4").Value = 3.65
$ws.Range("AA4").Value = 17.12
$ws.Range("AB4").Value = 20478.24
$ws.Range("AC4").Value = 41928
$ws.Range("AD4").Value = 17.51
$ws.Range("AE4").Value = 1050706
$ws.Range("AF4").Value = 0.7
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 0.14
$ws.Range("AI4").Value = 2.3
$ws.Range("AJ4").Value = 720000
$ws.Range("D5").Value = 11670
$ws.Range("E5").Value = 51
$ws.Range("F5").Value = 51
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = 50
$ws.Range("K5").Value = 10357
$ws.Range("L5").Value = 1407
$ws.Range("M5").Value = 8950
$ws.Range("N5").Value = 8950
$ws.Range("P5").Value = 44
$ws.Range("Q5").Value = 120
$ws.Range("R5").Value = -241
$ws.Range("S5").Value = -9
$ws.Range("T5").Value = 529
$ws.Range("U5").Value = -408
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0.43
$ws.Range("X5").Value = 0.43
$ws.Range("Y5").Value = 0.5600000000000001
$ws.Range("Z5").Value = 0.48
$ws.Range("AA5").Value = 15.72
$ws.Range("AB5").Value = 20608.54
$ws.Range("AC5").Value = 5663
$ws.Range("AD5").Value = 123.09
$ws.Range("AE5").Value = 1057421
$ws.Range("AF5").Value = 0.66
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 0.14
$ws.Range("AI5").Value = 17.02
$ws.Range("AJ5").Value = 720000
$ws.Range("D6").Value = 10797
$ws.Range("E6").Value = 86
$ws.Range("F6").Value = 86
$ws.Range("G6").Value = 77
$ws.Range("H6").Value = 20
$ws.Range("I6").Value = 20
$ws.Range("K6").Value = 10527
$ws.Range("L6").Value = 1638
$ws.Range("M6").Value = 8888
$ws.Range("N6").Value = 8888
$ws.Range("P6").Value = 44
$ws.Range("Q6").Value = 426
$ws.Range("R6").Value = -458
$ws.Range("S6").Value = -9
$ws.Range("T6").Value = 729
$ws.Range("U6").Value = -304
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0.8
$ws.Range("X6").Value = 0.19
$ws.Range("Y6").Value = 0.23
$ws.Range("Z6").Value = 0.19
$ws.Range("AA6").Value = 18.43
$ws.Range("AB6").Value = 20469.81
$ws.Range("AC6").Value = 2273
$ws.Range("AD6").Value = 274.54
$ws.Range("AE6").Value = 1050154
$ws.Range("AF6").Value = 0.59
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 0.16
$ws.Range("AI6").Value = 42.41
$ws.Range("AJ6").Value = 720000

# --- Remove cells that no longer exist (rows 2-5: J, O, and V for rows 2-4) ---
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# --- Rows 7-9: remove all data cells except A, B, C ---
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
